$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New meeting row for 10/10/2023 - copy date formatting from the row above
# then overwrite with the new date + attendance values.
$ws.Range("A3").Copy($ws.Range("A6"))
$ws.Range("A6").Value = (Get-Date -Year 2023 -Month 10 -Day 10).Date

$ws.Range("B6").Value = "Yes"
$ws.Range("C6").Value = "Yes"
$ws.Range("D6").Value = "Yes"
$ws.Range("E6").Value = "No"

$ws.Range("G11").Select()
